$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (shifts D:K -> E:L), preserving data.
$ws.Columns("D").Insert()

# The newly inserted column D has no formatting; copy formats from column E
# (which holds what used to be column D) so number/date formats match.
$ws.Range("E5:E102").Copy()
$ws.Range("D5:D102").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Populate the new column D with the new period's data.
$ws.Range("D7").Value2 = 43404
$ws.Range("D8").Value2 = 11100
$ws.Range("D9").Value2 = 9600
$ws.Range("D10").Value2 = 1600
$ws.Range("D12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("D15").Value2 = 0
$ws.Range("D17").Value2 = 10700
$ws.Range("D18").Value2 = 400
$ws.Range("D20").Value2 = 0
$ws.Range("D21").Value2 = 600
$ws.Range("D22").Value2 = "NA"
$ws.Range("D23").Value2 = 400
$ws.Range("D24").Value2 = 100
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 300
$ws.Range("D27").Value2 = 300
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = 0
$ws.Range("D33").Value2 = 300
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 300
$ws.Range("D38").Value2 = 43404
$ws.Range("D41").Value2 = 4700
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 600
$ws.Range("D44").Value2 = 4300
$ws.Range("D45").Value2 = 200
$ws.Range("D46").Value2 = 9700
$ws.Range("D47").Value2 = 0
$ws.Range("D48").Value2 = 1300
$ws.Range("D49").Value2 = 0
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 100
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 11000
$ws.Range("D57").Value2 = 600
$ws.Range("D58").Value2 = "NA"
$ws.Range("D59").Value2 = 1300
$ws.Range("D60").Value2 = 1900
$ws.Range("D61").Value2 = 0
$ws.Range("D62").Value2 = "NA"
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 1900
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 0
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = 3400
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 9200
$ws.Range("D77").Value2 = 0
$ws.Range("D80").Value2 = 43404
$ws.Range("D81").Value2 = 300
$ws.Range("D83").Value2 = 200
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 300
$ws.Range("D91").Value2 = -100
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -100
$ws.Range("D96").Value2 = 0
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = "NA"
$ws.Range("D101").Value2 = 0
$ws.Range("D102").Value2 = 200

# Special case: row 29 ("Discontinued Operations") - historical zero values
# for the shifted E:J columns become "NA" text instead of staying 0.
$ws.Range("E29:J29").Value2 = "NA"
